$wb = $excel.ActiveWorkbook

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 686.1429000000001
$ws.Range("I103").Value = 700.6
$ws.Range("J103").Value = 650
$ws.Range("K103").Value = 2101.8
$ws.Range("L103").Value = 1950
$ws.Range("M103").Value = -1515.8
$ws.Range("N103").Value = -3122

# ALC row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1781
$ws.Range("J107").Value = 1951.1538
$ws.Range("L107").Value = 1951.1538
$ws.Range("N107").Value = -5791.1538

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1872.6222
$ws.Range("I112").Value = 783.3333
$ws.Range("J112").Value = 2040.2051
$ws.Range("K112").Value = 2349.9999
$ws.Range("L112").Value = 6120.615299999999
$ws.Range("M112").Value = -1241.9999
$ws.Range("N112").Value = -8336.615299999999

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5381632.5
$ws.Range("I132").Value = 7410991.5
$ws.Range("J132").Value = 9800.647000000001
$ws.Range("K132").Value = 22232974.5
$ws.Range("L132").Value = 29401.941
$ws.Range("M132").Value = -22230444.5
$ws.Range("N132").Value = -34461.94100000001

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1050.8903
$ws.Range("I137").Value = 862.1556
$ws.Range("J137").Value = 1280.4324
$ws.Range("K137").Value = 2586.4668
$ws.Range("L137").Value = 3841.2972
$ws.Range("M137").Value = -36.46680000000015
$ws.Range("N137").Value = -8941.297200000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 551619.3
$ws.Range("I138").Value = 852.13513
$ws.Range("J138").Value = 1036818.94
$ws.Range("K138").Value = 2556.40539
$ws.Range("L138").Value = 3110456.82
$ws.Range("M138").Value = 2583.59461
$ws.Range("N138").Value = -3120736.82

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 19608780
$ws.Range("I61").Value = 21739924
$ws.Range("K61").Value = 21739924
$ws.Range("M61").Value = -21739712

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1001.55817
$ws.Range("J74").Value = 2475.4285
$ws.Range("L74").Value = 2475.4285
$ws.Range("N74").Value = -4223.4285

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1001.55817
$ws.Range("J77").Value = 2475.4285
$ws.Range("L77").Value = 12377.1425
$ws.Range("N77").Value = -21113.1425

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1631.0454
$ws.Range("I110").Value = 1226.3572
$ws.Range("K110").Value = 1226.3572
$ws.Range("M110").Value = 818.6428000000001

# ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 39799.5
$ws.Range("J125").Value = 39799.5
$ws.Range("L125").Value = 39799.5
$ws.Range("N125").Value = -49639.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2351.027
$ws.Range("I132").Value = 2549
$ws.Range("K132").Value = 7647
$ws.Range("M132").Value = -5117

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 19608780
$ws.Range("I136").Value = 21739924
$ws.Range("K136").Value = 65219772
$ws.Range("M136").Value = -65217222

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 848.21875
$ws.Range("I58").Value = 759.5417
$ws.Range("J58").Value = 1114.25
$ws.Range("K58").Value = 759.5417
$ws.Range("L58").Value = 1114.25
$ws.Range("M58").Value = -556.5417
$ws.Range("N58").Value = -1520.25

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 14961.308
$ws.Range("I60").Value = 3250
$ws.Range("J60").Value = 24999.572
$ws.Range("K60").Value = 3250
$ws.Range("L60").Value = 24999.572
$ws.Range("M60").Value = -2739
$ws.Range("N60").Value = -26021.572

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 29900
$ws.Range("J95").Value = 29900
$ws.Range("L95").Value = 29900
$ws.Range("N95").Value = -35392

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1229.2858
$ws.Range("I122").Value = 1107.2858
$ws.Range("K122").Value = 3321.8574
$ws.Range("M122").Value = -871.8574000000003

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 848.21875
$ws.Range("I136").Value = 759.5417
$ws.Range("J136").Value = 1114.25
$ws.Range("K136").Value = 2278.6251
$ws.Range("L136").Value = 3342.75
$ws.Range("M136").Value = 271.3748999999998
$ws.Range("N136").Value = -8442.75

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2179.0667
$ws.Range("J34").Value = 2599.4546
$ws.Range("L34").Value = 7798.3638
$ws.Range("N34").Value = -7966.3638

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3029.2942
$ws.Range("I39").Value = 3800
$ws.Range("J39").Value = 2864.1428
$ws.Range("K39").Value = 11400
$ws.Range("L39").Value = 8592.428400000001
$ws.Range("M39").Value = -11106
$ws.Range("N39").Value = -9180.428400000001

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2958.8
$ws.Range("I55").Value = 1004
$ws.Range("J55").Value = 3447.5
$ws.Range("K55").Value = 3012
$ws.Range("L55").Value = 10342.5
$ws.Range("M55").Value = -2835
$ws.Range("N55").Value = -10696.5

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15152703
$ws.Range("J131").Value = 1252.3442
$ws.Range("L131").Value = 3757.0326
$ws.Range("N131").Value = -13837.0326

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1482.825
$ws.Range("I132").Value = 1136.8485
$ws.Range("J132").Value = 3113.8572
$ws.Range("K132").Value = 3410.5455
$ws.Range("L132").Value = 9341.571599999999
$ws.Range("M132").Value = -880.5455000000002
$ws.Range("N132").Value = -14401.5716

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6699.8
$ws.Range("I46").Value = 6999.3335
$ws.Range("J46").Value = 6571.4287
$ws.Range("K46").Value = 6999.3335
$ws.Range("L46").Value = 6571.4287
$ws.Range("M46").Value = -6811.3335
$ws.Range("N46").Value = -6947.4287

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1819.6
$ws.Range("I61").Value = 1433
$ws.Range("J61").Value = 2399.5
$ws.Range("K61").Value = 1433
$ws.Range("L61").Value = 2399.5
$ws.Range("M61").Value = -1231
$ws.Range("N61").Value = -2803.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 834.3333
$ws.Range("J93").Value = 751.5
$ws.Range("L93").Value = 751.5
$ws.Range("N93").Value = -3247.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1819.6
$ws.Range("I113").Value = 1433
$ws.Range("J113").Value = 2399.5
$ws.Range("K113").Value = 1433
$ws.Range("L113").Value = 2399.5
$ws.Range("M113").Value = 737
$ws.Range("N113").Value = -6739.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 20838270
$ws.Range("I122").Value = 31251480
$ws.Range("J122").Value = 11850
$ws.Range("K122").Value = 93754440
$ws.Range("L122").Value = 35550
$ws.Range("M122").Value = -93751990
$ws.Range("N122").Value = -40450

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 23548.11
$ws.Range("I132").Value = 1216.8
$ws.Range("J132").Value = 50133
$ws.Range("K132").Value = 3650.4
$ws.Range("L132").Value = 150399
$ws.Range("M132").Value = -1120.4
$ws.Range("N132").Value = -155459

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1153.1724
$ws.Range("I136").Value = 1071.9259
$ws.Range("K136").Value = 3215.7777
$ws.Range("M136").Value = -665.7776999999996

# WVR row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17000
$ws.Range("J64").Value = 17000
$ws.Range("L64").Value = 17000
$ws.Range("N64").Value = -17496

# WVR row 67
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 17000
$ws.Range("J67").Value = 17000
$ws.Range("L67").Value = 17000
$ws.Range("N67").Value = -18716

# WVR row 127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 74600
$ws.Range("J127").Value = 74600
$ws.Range("L127").Value = 74600
$ws.Range("N127").Value = -84520

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2941.9167
$ws.Range("I132").Value = 3224.4814
$ws.Range("J132").Value = 2094.2222
$ws.Range("K132").Value = 9673.4442
$ws.Range("L132").Value = 6282.6666
$ws.Range("M132").Value = -7143.4442
$ws.Range("N132").Value = -11342.6666

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 596.0732
$ws.Range("I136").Value = 423.32257
$ws.Range("J136").Value = 1131.6
$ws.Range("K136").Value = 1269.96771
$ws.Range("L136").Value = 3394.8
$ws.Range("M136").Value = 1280.03229
$ws.Range("N136").Value = -8494.799999999999
